{"js": "// Update the three-digit-by-one-digit multiplication answers in the\n// worksheet table. Each \"old\" expression appears exactly once in the\n// document, so we can safely search-and-replace each pair in turn.\nconst replacements = [\n  [\"321\u00d72=642\", \"475\u00d75=2375\"],\n  [\"511\u00d75=2555\", \"959\u00d79=8631\"],\n  [\"340\u00d75=1700\", \"904\u00d76=5424\"],\n  [\"938\u00d73=2814\", \"394\u00d79=3546\"],\n  [\"944\u00d78=7552\", \"690\u00d74=2760\"],\n  [\"782\u00d76=4692\", \"591\u00d79=5319\"],\n  [\"899\u00d78=7192\", \"311\u00d76=1866\"],\n  [\"966\u00d77=6762\", \"152\u00d77=1064\"],\n  [\"907\u00d75=4535\", \"942\u00d79=8478\"],\n  [\"655\u00d75=3275\", \"756\u00d74=3024\"],\n  [\"537\u00d77=3759\", \"571\u00d77=3997\"],\n  [\"466\u00d74=1864\", \"220\u00d72=440\"],\n  [\"426\u00d75=2130\", \"544\u00d72=1088\"],\n  [\"401\u00d79=3609\", \"283\u00d78=2264\"],\n  [\"201\u00d79=1809\", \"784\u00d73=2352\"],\n  [\"559\u00d77=3913\", \"297\u00d72=594\"],\n  [\"940\u00d78=7520\", \"659\u00d72=1318\"],\n  [\"322\u00d75=1610\", \"999\u00d75=4995\"],\n  [\"623\u00d79=5607\", \"435\u00d76=2610\"],\n  [\"618\u00d77=4326\", \"205\u00d76=1230\"],\n  [\"296\u00d79=2664\", \"804\u00d76=4824\"],\n  [\"885\u00d77=6195\", \"528\u00d76=3168\"],\n  [\"114\u00d77=798\", \"292\u00d74=1168\"],\n  [\"965\u00d78=7720\", \"128\u00d76=768\"],\n  [\"277\u00d79=2493\", \"164\u00d79=1476\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find expected text \"${oldText}\" to replace.`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-by-one-digit multiplication answers in the\n# worksheet table. Each \"old\" expression occurs exactly once in the\n# document, so a simple Find/Replace per pair is sufficient and safe.\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"321\u00d72=642\", \"475\u00d75=2375\"),\n    @(\"511\u00d75=2555\", \"959\u00d79=8631\"),\n    @(\"340\u00d75=1700\", \"904\u00d76=5424\"),\n    @(\"938\u00d73=2814\", \"394\u00d79=3546\"),\n    @(\"944\u00d78=7552\", \"690\u00d74=2760\"),\n    @(\"782\u00d76=4692\", \"591\u00d79=5319\"),\n    @(\"899\u00d78=7192\", \"311\u00d76=1866\"),\n    @(\"966\u00d77=6762\", \"152\u00d77=1064\"),\n    @(\"907\u00d75=4535\", \"942\u00d79=8478\"),\n    @(\"655\u00d75=3275\", \"756\u00d74=3024\"),\n    @(\"537\u00d77=3759\", \"571\u00d77=3997\"),\n    @(\"466\u00d74=1864\", \"220\u00d72=440\"),\n    @(\"426\u00d75=2130\", \"544\u00d72=1088\"),\n    @(\"401\u00d79=3609\", \"283\u00d78=2264\"),\n    @(\"201\u00d79=1809\", \"784\u00d73=2352\"),\n    @(\"559\u00d77=3913\", \"297\u00d72=594\"),\n    @(\"940\u00d78=7520\", \"659\u00d72=1318\"),\n    @(\"322\u00d75=1610\", \"999\u00d75=4995\"),\n    @(\"623\u00d79=5607\", \"435\u00d76=2610\"),\n    @(\"618\u00d77=4326\", \"205\u00d76=1230\"),\n    @(\"296\u00d79=2664\", \"804\u00d76=4824\"),\n    @(\"885\u00d77=6195\", \"528\u00d76=3168\"),\n    @(\"114\u00d77=798\", \"292\u00d74=1168\"),\n    @(\"965\u00d78=7720\", \"128\u00d76=768\"),\n    @(\"277\u00d79=2493\", \"164\u00d79=1476\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne) | Out-Null\n}\n"}
